$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserDataSet")

# --- Clear stray "String.Empty" placeholder values that were mistakenly
# --- data-bound into the sheet, and align their formatting with the rest
# --- of the hyperlink-styled columns (style index 2) where applicable.

# Cells whose style must become the "hyperlink-font / centered" style
# (same as B2) in addition to clearing their value.
$styleSourceCell = $ws.Range("B2")
$restyleCells = @("C2","D4","E4","E7","C11","E11","C12","D12","E12","C13","E13","C14","E14","C15","E15","C16","E16")

$styleSourceCell.Copy()
foreach ($addr in $restyleCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

foreach ($addr in $restyleCells) {
    $ws.Range($addr).Value = ""
}

# Cells that just need their placeholder value cleared, keeping their
# existing style untouched.
$clearOnlyCells = @("B3","B11")
foreach ($addr in $clearOnlyCells) {
    $ws.Range($addr).Value = ""
}

# --- Reset the sheet view: drop the scrolled-right "topLeftCell" and move
# --- the active selection to B11 instead of E7.
$ws.Activate() | Out-Null
$ws.Range("B11").Select() | Out-Null
